$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1828.5488
$ws.Range("I15").Value = 1828.5488
$ws.Range("K15").Value = 5485.6464
$ws.Range("M15").Value = -5316.6464

$ws.Range("H40").Value = 1777.2727
$ws.Range("I40").Value = 1950
$ws.Range("J40").Value = 1678.5714
$ws.Range("K40").Value = 1950
$ws.Range("L40").Value = 1678.5714
$ws.Range("M40").Value = -1775
$ws.Range("N40").Value = -2028.5714

$ws.Range("H55").Value = 113.15385
$ws.Range("I55").Value = 68.5
$ws.Range("J55").Value = 151.42857
$ws.Range("K55").Value = 68.5
$ws.Range("L55").Value = 151.42857
$ws.Range("M55").Value = 145.5
$ws.Range("N55").Value = -579.42857

$ws.Range("H98").Value = 2821
$ws.Range("I98").Value = 1371.25
$ws.Range("J98").Value = 5720.5
$ws.Range("K98").Value = 1371.25
$ws.Range("L98").Value = 5720.5
$ws.Range("M98").Value = 126.75
$ws.Range("N98").Value = -8716.5

$ws.Range("H122").Value = 2821
$ws.Range("I122").Value = 1371.25
$ws.Range("J122").Value = 5720.5
$ws.Range("K122").Value = 4113.75
$ws.Range("L122").Value = 17161.5
$ws.Range("M122").Value = -1663.75
$ws.Range("N122").Value = -22061.5

$ws.Range("H137").Value = 1370.7858
$ws.Range("I137").Value = 1145.4706
$ws.Range("J137").Value = 1719
$ws.Range("K137").Value = 3436.4118
$ws.Range("L137").Value = 5157
$ws.Range("M137").Value = -886.4118000000003
$ws.Range("N137").Value = -10257

$ws.Range("H138").Value = 536275.6
$ws.Range("I138").Value = 3176.4167
$ws.Range("J138").Value = 1028367.1
$ws.Range("K138").Value = 9529.250100000001
$ws.Range("L138").Value = 3085101.3
$ws.Range("M138").Value = -4389.250100000001
$ws.Range("N138").Value = -3095381.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H61").Value = 1978.3214
$ws.Range("I61").Value = 1681.2222
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 1681.2222
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -1469.2222
$ws.Range("N61").Value = -10424

$ws.Range("H97").Value = 2632.825
$ws.Range("I97").Value = 1620
$ws.Range("J97").Value = 3645.65
$ws.Range("K97").Value = 1620
$ws.Range("L97").Value = 3645.65
$ws.Range("M97").Value = -1124
$ws.Range("N97").Value = -4637.65

$ws.Range("H132").Value = 1926.3077
$ws.Range("I132").Value = 1746.2391
$ws.Range("J132").Value = 3306.8333
$ws.Range("K132").Value = 5238.7173
$ws.Range("L132").Value = 9920.499899999999
$ws.Range("M132").Value = -2708.7173
$ws.Range("N132").Value = -14980.4999

$ws.Range("H136").Value = 1978.3214
$ws.Range("I136").Value = 1681.2222
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 5043.6666
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -2493.6666
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7577702.5
$ws.Range("I20").Value = 15153758
$ws.Range("J20").Value = 1646.5454
$ws.Range("K20").Value = 15153758
$ws.Range("L20").Value = 1646.5454
$ws.Range("M20").Value = -15153511
$ws.Range("N20").Value = -2140.5454

$ws.Range("H22").Value = 1350
$ws.Range("I22").Value = 1350
$ws.Range("K22").Value = 1350
$ws.Range("M22").Value = -1177

$ws.Range("H134").Value = 2303.8276
$ws.Range("I134").Value = 2366.5557
$ws.Range("J134").Value = 1457
$ws.Range("K134").Value = 7099.6671
$ws.Range("L134").Value = 4371
$ws.Range("M134").Value = -4564.6671
$ws.Range("N134").Value = -9441

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1696.0952
$ws.Range("I58").Value = 1571.7333
$ws.Range("J58").Value = 2007
$ws.Range("K58").Value = 1571.7333
$ws.Range("L58").Value = 2007
$ws.Range("M58").Value = -1368.7333
$ws.Range("N58").Value = -2413

$ws.Range("H132").Value = 3682.111
$ws.Range("I132").Value = 3098.5715
$ws.Range("J132").Value = 5724.5
$ws.Range("K132").Value = 9295.7145
$ws.Range("L132").Value = 17173.5
$ws.Range("M132").Value = -6765.7145
$ws.Range("N132").Value = -22233.5

$ws.Range("H134").Value = 24556.715
$ws.Range("I134").Value = 24556.715
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 73670.145
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -71135.145

$ws.Range("H136").Value = 1696.0952
$ws.Range("I136").Value = 1571.7333
$ws.Range("J136").Value = 2007
$ws.Range("K136").Value = 4715.199900000001
$ws.Range("L136").Value = 6021
$ws.Range("M136").Value = -2165.199900000001
$ws.Range("N136").Value = -11121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 926.2
$ws.Range("I68").Value = 721.7963
$ws.Range("J68").Value = 1166.1522
$ws.Range("K68").Value = 2165.3889
$ws.Range("L68").Value = 3498.4566
$ws.Range("M68").Value = -1354.3889
$ws.Range("N68").Value = -5120.4566

$ws.Range("H71").Value = 926.2
$ws.Range("I71").Value = 721.7963
$ws.Range("J71").Value = 1166.1522
$ws.Range("K71").Value = 6496.1667
$ws.Range("L71").Value = 10495.3698
$ws.Range("M71").Value = -2440.1667
$ws.Range("N71").Value = -18607.3698

$ws.Range("H129").Value = 46855.59
$ws.Range("J129").Value = 1803.5834
$ws.Range("L129").Value = 5410.7502
$ws.Range("N129").Value = -15410.7502

$ws.Range("H131").Value = 1185.0571
$ws.Range("I131").Value = 430.66666
$ws.Range("J131").Value = 1578.6522
$ws.Range("K131").Value = 1291.99998
$ws.Range("L131").Value = 4735.9566
$ws.Range("M131").Value = 3748.00002
$ws.Range("N131").Value = -14815.9566

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 29600.404
$ws.Range("I70").Value = 54257.523
$ws.Range("J70").Value = 4943.2856
$ws.Range("K70").Value = 54257.523
$ws.Range("L70").Value = 4943.2856
$ws.Range("M70").Value = -53987.523
$ws.Range("N70").Value = -5483.2856

$ws.Range("H73").Value = 29600.404
$ws.Range("I73").Value = 54257.523
$ws.Range("J73").Value = 4943.2856
$ws.Range("K73").Value = 54257.523
$ws.Range("L73").Value = 4943.2856
$ws.Range("M73").Value = -53321.523
$ws.Range("N73").Value = -6815.2856

$ws.Range("H97").Value = 1593
$ws.Range("I97").Value = 1593
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1593
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -1097

$ws.Range("H102").Value = 1868.1177
$ws.Range("I102").Value = 1588
$ws.Range("J102").Value = 2778.5
$ws.Range("K102").Value = 1588
$ws.Range("L102").Value = 2778.5
$ws.Range("M102").Value = 34
$ws.Range("N102").Value = -6022.5

$ws.Range("H113").Value = 1544.4375
$ws.Range("I113").Value = 1594.0667
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 1594.0667
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 575.9332999999999
$ws.Range("N113").Value = -5140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1034.5883
$ws.Range("I22").Value = 829.75
$ws.Range("J22").Value = 1216.6666
$ws.Range("K22").Value = 829.75
$ws.Range("L22").Value = 1216.6666
$ws.Range("M22").Value = -534.75
$ws.Range("N22").Value = -1806.6666

$ws.Range("H27").Value = 1034.5883
$ws.Range("I27").Value = 829.75
$ws.Range("J27").Value = 1216.6666
$ws.Range("K27").Value = 829.75
$ws.Range("L27").Value = 1216.6666
$ws.Range("M27").Value = -722.75
$ws.Range("N27").Value = -1430.6666

$ws.Range("H93").Value = 1913.3334
$ws.Range("I93").Value = 1948
$ws.Range("J93").Value = 1867.1111
$ws.Range("K93").Value = 1948
$ws.Range("L93").Value = 1867.1111
$ws.Range("M93").Value = -700
$ws.Range("N93").Value = -4363.1111

$ws.Range("H122").Value = 3034.087
$ws.Range("I122").Value = 2260.4
$ws.Range("J122").Value = 3629.2307
$ws.Range("K122").Value = 6781.200000000001
$ws.Range("L122").Value = 10887.6921
$ws.Range("M122").Value = -4331.200000000001
$ws.Range("N122").Value = -15787.6921

$ws.Range("H132").Value = 2098.8
$ws.Range("I132").Value = 1434.1428
$ws.Range("J132").Value = 3649.6667
$ws.Range("K132").Value = 4302.428400000001
$ws.Range("L132").Value = 10949.0001
$ws.Range("M132").Value = -1772.428400000001
$ws.Range("N132").Value = -16009.0001

$ws.Range("H136").Value = 1991.5555
$ws.Range("I136").Value = 1859.1428
$ws.Range("J136").Value = 2455
$ws.Range("K136").Value = 5577.428400000001
$ws.Range("L136").Value = 7365
$ws.Range("M136").Value = -3027.428400000001
$ws.Range("N136").Value = -12465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2125
$ws.Range("I96").Value = 1966.6666
$ws.Range("J96").Value = 2600
$ws.Range("K96").Value = 1966.6666
$ws.Range("L96").Value = 2600
$ws.Range("M96").Value = -593.6666
$ws.Range("N96").Value = -5346

$ws.Range("H122").Value = 11071.48
$ws.Range("I122").Value = 18661.385
$ws.Range("J122").Value = 2849.0833
$ws.Range("K122").Value = 55984.155
$ws.Range("L122").Value = 8547.249899999999
$ws.Range("M122").Value = -53534.155
$ws.Range("N122").Value = -13447.2499

$ws.Range("H132").Value = 11333.667
$ws.Range("I132").Value = 14501
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 43503
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -40973
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 1310.2128
$ws.Range("I136").Value = 1043.2894
$ws.Range("J136").Value = 2437.2222
$ws.Range("K136").Value = 3129.8682
$ws.Range("L136").Value = 7311.6666
$ws.Range("M136").Value = -579.8681999999999
$ws.Range("N136").Value = -12411.6666
